$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1 dashboard")
$ws2 = $wb.Worksheets.Item("2 signup")

$sheet1Labels = @(
    "userDatabase",
    "logout",
    "message",
    "newLanguage",
    "newPage",
    "newLabel",
    "updateLabel",
    "fullTemplate",
    "downloadTemplates",
    "completeUpload"
)

$sheet2Labels = @(
    "signUp",
    "userName",
    "password",
    "login",
    "firstName",
    "lastName",
    "email",
    "userDatabase"
)

for ($i = 0; $i -lt $sheet1Labels.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $sheet1Labels[$i]
}

for ($i = 0; $i -lt $sheet2Labels.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $sheet2Labels[$i]
}

$ws1.Range("C10:C500").Validation.Delete()
$ws1.Range("C2:C500").Validation.Delete()
$ws2.Range("C10:C500").Validation.Delete()
$ws2.Range("C2:C500").Validation.Delete()
